$wb = $excel.ActiveWorkbook

# --- 1. Remove the empty B5 / B10 cells on the "ODI Batting" sheet -------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("B5").ClearContents()
$battingSheet.Range("B10").ClearContents()

# --- 2. Add the new "ODI Batting Extra" worksheet after "ODI Bowling" ----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Make every cell we touch hold plain text (matches the inline-string data in
# the source file) unless it is explicitly meant to be numeric (column B,
# BATTING_POSITION, stays in its default General number format for the rows
# that carry a real number; rows 2 and 10 leave BATTING_POSITION blank so we
# format those two cells as text too, to keep them present-but-empty).
$newSheet.Range("A1:A10").NumberFormat = "@"
$newSheet.Range("C1:F10").NumberFormat = "@"
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B10").NumberFormat = "@"

# Header row
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Reuse the same header style (bold, bordered, centered) already used by the
# other sheets' header rows.
$headerStyleSource = $wb.Worksheets.Item("Player Info").Range("A1")
$headerStyleSource.Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)

# Data rows (A, C, D, E, F as text; B as a real number when present)
$newSheet.Range("A2").Value = "4258"
$newSheet.Range("B2").Value = ""
$newSheet.Range("C2").Value = ""
$newSheet.Range("D2").Value = ""
$newSheet.Range("E2").Value = ""
$newSheet.Range("F2").Value = "NO"

$newSheet.Range("A3").Value = "4268"
$newSheet.Range("B3").Value = 6
$newSheet.Range("C3").Value = "5"
$newSheet.Range("D3").Value = "6"
$newSheet.Range("E3").Value = "23.40%"
$newSheet.Range("F3").Value = "YES"

$newSheet.Range("A4").Value = "4270"
$newSheet.Range("B4").Value = 6
$newSheet.Range("C4").Value = "2"
$newSheet.Range("D4").Value = "1"
$newSheet.Range("E4").Value = "7.35%"
$newSheet.Range("F4").Value = "NO"

$newSheet.Range("A5").Value = "4398"
$newSheet.Range("B5").Value = 6
$newSheet.Range("C5").Value = ""
$newSheet.Range("D5").Value = ""
$newSheet.Range("E5").Value = ""
$newSheet.Range("F5").Value = "NO"

$newSheet.Range("A6").Value = "4399"
$newSheet.Range("B6").Value = 6
$newSheet.Range("C6").Value = "0"
$newSheet.Range("D6").Value = "0"
$newSheet.Range("E6").Value = "4.28%"
$newSheet.Range("F6").Value = "NO"

$newSheet.Range("A7").Value = "4400"
$newSheet.Range("B7").Value = 7
$newSheet.Range("C7").Value = "0"
$newSheet.Range("D7").Value = "0"
$newSheet.Range("E7").Value = "1.40%"
$newSheet.Range("F7").Value = "NO"

$newSheet.Range("A8").Value = "4483"
$newSheet.Range("B8").Value = 6
$newSheet.Range("C8").Value = "2"
$newSheet.Range("D8").Value = "2"
$newSheet.Range("E8").Value = "19.44%"
$newSheet.Range("F8").Value = "NO"

$newSheet.Range("A9").Value = "4484"
$newSheet.Range("B9").Value = 6
$newSheet.Range("C9").Value = "0"
$newSheet.Range("D9").Value = "0"
$newSheet.Range("E9").Value = "0.53%"
$newSheet.Range("F9").Value = "NO"

$newSheet.Range("A10").Value = "4486"
$newSheet.Range("B10").Value = ""
$newSheet.Range("C10").Value = ""
$newSheet.Range("D10").Value = ""
$newSheet.Range("E10").Value = ""
$newSheet.Range("F10").Value = "NO"
